# Bug fix for a non-text file specified in a rule
# Insert a new test-case row (row 38) on the CheckList sheet, pushing the
# existing rows 38-67 down to 39-68, and set the new row's contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CheckList")

# Insert a new row above the current row 38 - this shifts rows 38:67 down to 39:68.
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new test case.
$ws.Cells.Item(38, 1).Value = "Functionality"
$ws.Cells.Item(38, 2).Value = "If opening a solution with a rule for an image file, then nothing is modified & an error is duly reported"
$ws.Cells.Item(38, 3).Value = "n/a"
$ws.Cells.Item(38, 4).Value = "n/a"
$ws.Cells.Item(38, 5).Value = "n/a"
$ws.Cells.Item(38, 6).Value = "n/a"
$ws.Cells.Item(38, 7).Value = "n/a"
$ws.Cells.Item(38, 8).Value = "to be tested"

# Match styling used by the other rows in this section: plain wrap-text for
# A/B, green "n/a" fill for C:G, yellow "TODO/to be tested" fill for H.
$ws.Range("A38").Style = $ws.Range("A37").Style
$ws.Range("B38").Style = $ws.Range("A37").Style
$ws.Range("C38:G38").Style = $ws.Range("C37:G37").Style
$ws.Range("H38").Style = $ws.Range("H48").Style

# This row's text wraps over two lines, like other long-text rows.
$ws.Rows.Item(38).RowHeight = 22.5

# Update the frozen-pane top-left cell and the active selection to match the
# view position recorded after the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("H59").Select()
